# Update the header cell C1 on every "Unit N" sheet from "example_usage" to "example",
# and restore each sheet's last-known selection/active cell as captured by the author.

$wb = $excel.ActiveWorkbook

# 1) Fix header text on every sheet (C1: "example_usage" -> "example")
foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "example"
}

# 2) Re-apply the per-sheet selection / active-cell state.
$wsUnit1 = $wb.Worksheets.Item("Unit 1")
$wsUnit1.Activate() | Out-Null
$wsUnit1.Range("E16").Select() | Out-Null

$wsUnit2 = $wb.Worksheets.Item("Unit 2")
$wsUnit2.Activate() | Out-Null
$wsUnit2.Range("A1:B1").Select() | Out-Null

$wsUnit5 = $wb.Worksheets.Item("Unit 5")
$wsUnit5.Activate() | Out-Null
$wsUnit5.Range("A1:B1").Select() | Out-Null

$wsUnit7 = $wb.Worksheets.Item("Unit 7")
$wsUnit7.Activate() | Out-Null
$wsUnit7.Range("D14").Select() | Out-Null

# Leave "Unit 1" as the active/visible tab, matching tabSelected="1" in the source.
$wsUnit1.Activate() | Out-Null
